$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new cells to keep their content as TEXT (matches t="inlineStr"/t="s" text
# cells in the target -- "80" stays the text "80", not the number 80).
$ws.Range("A2:C4").NumberFormat = "@"

# New data rows published below the header row
$ws.Range("A2").Value = "Honduras"
$ws.Range("B2").Value = "80"
$ws.Range("C2").Value = "8703.10"

$ws.Range("A3").Value = "Guatemala"
$ws.Range("B3").Value = "80"
$ws.Range("C3").Value = "8344.90"

$ws.Range("A4").Value = "Costa Rica"
$ws.Range("B4").Value = "80"
$ws.Range("C4").Value = "9503.10"

# Ensure the new cells use the default (unstyled) cell style, same as the target XML
# which has no explicit style override for the data rows (s="0").
$ws.Range("A2:C4").Style = "Normal"

# Sheet view stays left-to-right (matches rightToLeft="false" in the target).
$excel.ActiveWindow.DisplayRightToLeft = $false
